$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("Home_Index_Title", "Startseite"),
    @("Home_Index_Fight", "Kampf"),
    @("Home_Index_XpCalc", "Erfahrung"),
    @("Home_Index_MyPokemon", "Mein Pokémon"),
    @("Home_Index_OpponentPokemon", "Gegnerisches Pokémon"),
    @("Home_Index_UsedAttack", "Eingesetzte Attacke"),
    @("Home_Index_Level", "Level"),
    @("Home_Index_Pokemon", "Pokémon"),
    @("Home_Index_SearchPokemon", "Pokémon suchen")
)

$startRow = 141
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $rows[$i][0]
    $ws.Cells.Item($r, 2).Value = $rows[$i][1]
}

$ws.Range("A149:B149").Select()
